$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet tab name / title to reflect new "through" date
$ws.Name = "Through 2021-09-25"

# Row 6 (April) - 2020 columns (Q,R,S)
$ws.Range("Q6").Value = 3
$ws.Range("R6").Value = 61
$ws.Range("S6").Value = 0.0469

# Row 11 (September) - update label text
$ws.Range("A11").Value = "September (through 09-25)"

# 2016 columns (E,F,G)
$ws.Range("F11").Value = 34
$ws.Range("G11").Value = 0.0556

# 2017 columns (H,I,J)
$ws.Range("I11").Value = 60
$ws.Range("J11").Value = 0.0625

# 2018 columns (K,L,M)
$ws.Range("L11").Value = 45
$ws.Range("M11").Value = 0.0816

# 2019 columns (N,O,P)
$ws.Range("O11").Value = 55
$ws.Range("P11").Value = 0.0984

# 2020 columns (Q,R,S)
$ws.Range("R11").Value = 95
$ws.Range("S11").Value = 0.0306

# 2021 columns (T,U,V)
$ws.Range("T11").Value = 2
$ws.Range("U11").Value = 150
$ws.Range("V11").Value = 0.0132

# Row 12 (Total)
# 2016 columns (E,F,G)
$ws.Range("F12").Value = 374
$ws.Range("G12").Value = 0.1031

# 2017 columns (H,I,J)
$ws.Range("I12").Value = 566
$ws.Range("J12").Value = 0.0797

# 2018 columns (K,L,M)
$ws.Range("L12").Value = 478
$ws.Range("M12").Value = 0.1132

# 2019 columns (N,O,P)
$ws.Range("O12").Value = 368
$ws.Range("P12").Value = 0.1024

# 2020 columns (Q,R,S)
$ws.Range("Q12").Value = 50
$ws.Range("R12").Value = 832
$ws.Range("S12").Value = 0.0567

# 2021 columns (T,U,V)
$ws.Range("T12").Value = 76
$ws.Range("U12").Value = 1146
$ws.Range("V12").Value = 0.0622
